# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E)
# for rows 2-51, per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.940.63' },
    @{ Cell = 'E2'; Value = '  +1.06%  ' },
    @{ Cell = 'D3'; Value = '1.877.46' },
    @{ Cell = 'E3'; Value = '  +0.36%  ' },
    @{ Cell = 'D4'; Value = '1.011' },
    @{ Cell = 'E4'; Value = '  +0.69%  ' },
    @{ Cell = 'D5'; Value = '335.97' },
    @{ Cell = 'E5'; Value = '  +1.23%  ' },
    @{ Cell = 'D6'; Value = '1.011' },
    @{ Cell = 'E6'; Value = '  +0.76%  ' },
    @{ Cell = 'D7'; Value = '0.4757' },
    @{ Cell = 'E7'; Value = '  +0.70%  ' },
    @{ Cell = 'D8'; Value = '0.3942' },
    @{ Cell = 'E8'; Value = '  +0.00%  ' },
    @{ Cell = 'D9'; Value = '46.96' },
    @{ Cell = 'E9'; Value = '  -2.15%  ' },
    @{ Cell = 'D10'; Value = '0.08000' },
    @{ Cell = 'E10'; Value = '  -0.60%  ' },
    @{ Cell = 'E11'; Value = '  -0.82%  ' },
    @{ Cell = 'E12'; Value = '  -0.68%  ' },
    @{ Cell = 'D13'; Value = '1.888.62' },
    @{ Cell = 'E13'; Value = '  -0.13%  ' },
    @{ Cell = 'D14'; Value = '6.038' },
    @{ Cell = 'E14'; Value = '  +1.46%  ' },
    @{ Cell = 'D15'; Value = '7.196' },
    @{ Cell = 'E15'; Value = '  +0.79%  ' },
    @{ Cell = 'E16'; Value = '  +0.73%  ' },
    @{ Cell = 'D17'; Value = '88.36' },
    @{ Cell = 'E17'; Value = '  +1.77%  ' },
    @{ Cell = 'D18'; Value = '0.06743' },
    @{ Cell = 'E18'; Value = '  +1.70%  ' },
    @{ Cell = 'E19'; Value = '  +0.34%  ' },
    @{ Cell = 'E20'; Value = '  -0.98%  ' },
    @{ Cell = 'E21'; Value = '  +0.60%  ' },
    @{ Cell = 'D22'; Value = '27.936.50' },
    @{ Cell = 'E22'; Value = '  +1.01%  ' },
    @{ Cell = 'D23'; Value = '5.500' },
    @{ Cell = 'E23'; Value = '  -0.07%  ' },
    @{ Cell = 'E24'; Value = '  +0.03%  ' },
    @{ Cell = 'D25'; Value = '2.342' },
    @{ Cell = 'E25'; Value = '  +1.53%  ' },
    @{ Cell = 'D26'; Value = '2.109.43' },
    @{ Cell = 'E26'; Value = '  -0.14%  ' },
    @{ Cell = 'D27'; Value = '158.45' },
    @{ Cell = 'E27'; Value = '  -0.21%  ' },
    @{ Cell = 'D28'; Value = '19.85' },
    @{ Cell = 'E28'; Value = '  -2.03%  ' },
    @{ Cell = 'D29'; Value = '2.099' },
    @{ Cell = 'E29'; Value = '  +0.12%  ' },
    @{ Cell = 'D30'; Value = '5.458' },
    @{ Cell = 'E30'; Value = '  -1.81%  ' },
    @{ Cell = 'D31'; Value = '121.35' },
    @{ Cell = 'E31'; Value = '  -0.76%  ' },
    @{ Cell = 'D32'; Value = '0.9734' },
    @{ Cell = 'E32'; Value = '  +0.48%  ' },
    @{ Cell = 'D33'; Value = '0.09526' },
    @{ Cell = 'E33'; Value = '  -0.08%  ' },
    @{ Cell = 'D34'; Value = '3.636' },
    @{ Cell = 'E34'; Value = '  +1.31%  ' },
    @{ Cell = 'D35'; Value = '5.332' },
    @{ Cell = 'E35'; Value = '  -0.10%  ' },
    @{ Cell = 'D36'; Value = '1.351' },
    @{ Cell = 'E36'; Value = '  -6.89%  ' },
    @{ Cell = 'D37'; Value = '0.06083' },
    @{ Cell = 'E37'; Value = '  -0.36%  ' },
    @{ Cell = 'D38'; Value = '0.02234' },
    @{ Cell = 'E38'; Value = '  -1.09%  ' },
    @{ Cell = 'E39'; Value = '  -1.84%  ' },
    @{ Cell = 'D40'; Value = '8.171' },
    @{ Cell = 'E40'; Value = '  +0.22%  ' },
    @{ Cell = 'D41'; Value = '1.010' },
    @{ Cell = 'E41'; Value = '  +0.70%  ' },
    @{ Cell = 'D42'; Value = '0.5971' },
    @{ Cell = 'E42'; Value = '  -1.06%  ' },
    @{ Cell = 'D43'; Value = '0.1893' },
    @{ Cell = 'E43'; Value = '  +0.01%  ' },
    @{ Cell = 'D44'; Value = '10.33' },
    @{ Cell = 'E44'; Value = '  +0.78%  ' },
    @{ Cell = 'D45'; Value = '1.266' },
    @{ Cell = 'E45'; Value = '  +0.77%  ' },
    @{ Cell = 'D46'; Value = '0.5667' },
    @{ Cell = 'E46'; Value = '  -0.81%  ' },
    @{ Cell = 'D47'; Value = '12.13' },
    @{ Cell = 'E47'; Value = '  -1.02%  ' },
    @{ Cell = 'E48'; Value = '  -0.73%  ' },
    @{ Cell = 'D49'; Value = '3.337' },
    @{ Cell = 'E49'; Value = '  -1.24%  ' },
    @{ Cell = 'D50'; Value = '0.06789' },
    @{ Cell = 'E50'; Value = '  -1.09%  ' },
    @{ Cell = 'D51'; Value = '112.23' },
    @{ Cell = 'E51'; Value = '  -1.95%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Value -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number (e.g. "1.011", "0.08000", "5.500") -- force
        # text storage first so Excel doesn't reinterpret/round it as a double,
        # then drop back to the default style so no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
